$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 519.3125
$ws.Range("I32").Value = 516.6667
$ws.Range("J32").Value = 520.9
$ws.Range("K32").Value = 516.6667
$ws.Range("L32").Value = 520.9
$ws.Range("M32").Value = -190.6667
$ws.Range("N32").Value = -1172.9
$ws.Range("H40").Value = 951.0417
$ws.Range("I40").Value = 773.2143
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 773.2143
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -598.2143
$ws.Range("N40").Value = -1550
$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 7999.5
$ws.Range("K51").Value = 7999.5
$ws.Range("M51").Value = -7515.5
$ws.Range("H64").Value = 4100
$ws.Range("J64").Value = 4500
$ws.Range("L64").Value = 4500
$ws.Range("N64").Value = -4996
$ws.Range("H67").Value = 4100
$ws.Range("J67").Value = 4500
$ws.Range("L67").Value = 4500
$ws.Range("N67").Value = -6216
$ws.Range("H111").Value = 3237.0715
$ws.Range("J111").Value = 1823
$ws.Range("L111").Value = 5469
$ws.Range("N111").Value = -11603
$ws.Range("H129").Value = 323543.28
$ws.Range("I129").Value = 266.66666
$ws.Range("J129").Value = 401129.7
$ws.Range("K129").Value = 799.9999799999999
$ws.Range("L129").Value = 1203389.1
$ws.Range("M129").Value = 4200.00002
$ws.Range("N129").Value = -1213389.1
$ws.Range("H132").Value = 21278742
$ws.Range("I132").Value = 22729462
$ws.Range("J132").Value = 1502
$ws.Range("K132").Value = 68188386
$ws.Range("L132").Value = 4506
$ws.Range("M132").Value = -68185856
$ws.Range("N132").Value = -9566
$ws.Range("H141").Value = 4119.5625
$ws.Range("I141").Value = 3727.5334
$ws.Range("K141").Value = 11182.6002
$ws.Range("M141").Value = -6002.600199999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1848.4375
$ws.Range("I2").Value = 1798.8182
$ws.Range("J2").Value = 1957.6
$ws.Range("K2").Value = 1798.8182
$ws.Range("L2").Value = 1957.6
$ws.Range("M2").Value = -1685.8182
$ws.Range("N2").Value = -2183.6
$ws.Range("H32").Value = 14239.818
$ws.Range("I32").Value = 9810.4375
$ws.Range("J32").Value = 26051.5
$ws.Range("K32").Value = 9810.4375
$ws.Range("L32").Value = 26051.5
$ws.Range("M32").Value = -9523.4375
$ws.Range("N32").Value = -26625.5
$ws.Range("H45").Value = 4928.8667
$ws.Range("I45").Value = 5133.8887
$ws.Range("K45").Value = 5133.8887
$ws.Range("M45").Value = -4756.8887
$ws.Range("H74").Value = 30304852
$ws.Range("I74").Value = 50000784
$ws.Range("K74").Value = 50000784
$ws.Range("M74").Value = -49999910
$ws.Range("H77").Value = 30304852
$ws.Range("I77").Value = 50000784
$ws.Range("K77").Value = 250003920
$ws.Range("M77").Value = -249999552
$ws.Range("H116").Value = 1848.4375
$ws.Range("I116").Value = 1798.8182
$ws.Range("J116").Value = 1957.6
$ws.Range("K116").Value = 1798.8182
$ws.Range("L116").Value = 1957.6
$ws.Range("M116").Value = 495.1818000000001
$ws.Range("N116").Value = -6545.6
$ws.Range("H125").Value = 34915
$ws.Range("J125").Value = 34915
$ws.Range("L125").Value = 34915
$ws.Range("N125").Value = -44755
$ws.Range("H132").Value = 10651784
$ws.Range("I132").Value = 13891342
$ws.Range("K132").Value = 41674026
$ws.Range("M132").Value = -41671496

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1848.4375
$ws.Range("I3").Value = 1798.8182
$ws.Range("J3").Value = 1957.6
$ws.Range("K3").Value = 1798.8182
$ws.Range("L3").Value = 1957.6
$ws.Range("M3").Value = -1684.8182
$ws.Range("N3").Value = -2185.6
$ws.Range("H86").Value = 2052.2354
$ws.Range("I86").Value = 1872.0667
$ws.Range("J86").Value = 3403.5
$ws.Range("K86").Value = 1872.0667
$ws.Range("L86").Value = 3403.5
$ws.Range("M86").Value = -749.0667000000001
$ws.Range("N86").Value = -5649.5
$ws.Range("H89").Value = 2052.2354
$ws.Range("I89").Value = 1872.0667
$ws.Range("J89").Value = 3403.5
$ws.Range("K89").Value = 9360.333500000001
$ws.Range("L89").Value = 17017.5
$ws.Range("M89").Value = -3744.333500000001
$ws.Range("N89").Value = -28249.5
$ws.Range("H134").Value = 3294.6296
$ws.Range("I134").Value = 3086.8635
$ws.Range("K134").Value = 9260.5905
$ws.Range("M134").Value = -6725.5905

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5749.3125
$ws.Range("I31").Value = 3097.9285
$ws.Range("J31").Value = 6841.0586
$ws.Range("K31").Value = 3097.9285
$ws.Range("L31").Value = 6841.0586
$ws.Range("M31").Value = -2802.9285
$ws.Range("N31").Value = -7431.0586
$ws.Range("H34").Value = 5749.3125
$ws.Range("I34").Value = 3097.9285
$ws.Range("J34").Value = 6841.0586
$ws.Range("K34").Value = 3097.9285
$ws.Range("L34").Value = 6841.0586
$ws.Range("M34").Value = -2895.9285
$ws.Range("N34").Value = -7245.0586
$ws.Range("H58").Value = 17895.227
$ws.Range("I58").Value = 1861.3334
$ws.Range("J58").Value = 28021.895
$ws.Range("K58").Value = 1861.3334
$ws.Range("L58").Value = 28021.895
$ws.Range("M58").Value = -1658.3334
$ws.Range("N58").Value = -28427.895
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H122").Value = 1523
$ws.Range("I122").Value = 1250.6666
$ws.Range("J122").Value = 1768.1
$ws.Range("K122").Value = 3751.9998
$ws.Range("L122").Value = 5304.299999999999
$ws.Range("M122").Value = -1301.9998
$ws.Range("N122").Value = -10204.3
$ws.Range("H132").Value = 40003840
$ws.Range("I132").Value = 50002850
$ws.Range("K132").Value = 150008550
$ws.Range("M132").Value = -150006020
$ws.Range("H134").Value = 43479200
$ws.Range("I134").Value = 45455436
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 136366308
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -136363773
$ws.Range("N134").Value = -11070
$ws.Range("H136").Value = 17895.227
$ws.Range("I136").Value = 1861.3334
$ws.Range("J136").Value = 28021.895
$ws.Range("K136").Value = 5584.0002
$ws.Range("L136").Value = 84065.685
$ws.Range("M136").Value = -3034.0002
$ws.Range("N136").Value = -89165.685

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()
$ws.Range("H122").Value = 1376.6897
$ws.Range("J122").Value = 1482.5385
$ws.Range("L122").Value = 13342.8465
$ws.Range("N122").Value = -18242.8465
$ws.Range("H131").Value = 714.99
$ws.Range("J131").Value = 752.0112
$ws.Range("L131").Value = 2256.0336
$ws.Range("N131").Value = -12336.0336
$ws.Range("H134").Value = 4182.7
$ws.Range("I134").Value = 3102.7273
$ws.Range("J134").Value = 5502.6665
$ws.Range("K134").Value = 9308.1819
$ws.Range("L134").Value = 16507.9995
$ws.Range("M134").Value = -4238.1819
$ws.Range("N134").Value = -26647.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 13448
$ws.Range("J94").Value = 13448
$ws.Range("L94").Value = 13448
$ws.Range("N94").Value = -14800
$ws.Range("H102").Value = 2862.0557
$ws.Range("I102").Value = 2531.4375
$ws.Range("J102").Value = 5507
$ws.Range("K102").Value = 2531.4375
$ws.Range("L102").Value = 5507
$ws.Range("M102").Value = -909.4375
$ws.Range("N102").Value = -8751
$ws.Range("H113").Value = 10502.333
$ws.Range("I113").Value = 11665.125
$ws.Range("K113").Value = 11665.125
$ws.Range("M113").Value = -9495.125
$ws.Range("H132").Value = 3113029.8
$ws.Range("I132").Value = 3851539.8
$ws.Range("K132").Value = 11554619.4
$ws.Range("M132").Value = -11552089.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1190.2
$ws.Range("I46").Value = 985.4286
$ws.Range("J46").Value = 1668
$ws.Range("K46").Value = 985.4286
$ws.Range("L46").Value = 1668
$ws.Range("M46").Value = -797.4286
$ws.Range("N46").Value = -2044
$ws.Range("H61").Value = 4657.5713
$ws.Range("I61").Value = 2443.1428
$ws.Range("K61").Value = 2443.1428
$ws.Range("M61").Value = -2241.1428
$ws.Range("H68").Value = 2450.3333
$ws.Range("J68").Value = 2934
$ws.Range("L68").Value = 2934
$ws.Range("N68").Value = -4432
$ws.Range("H71").Value = 2450.3333
$ws.Range("J71").Value = 2934
$ws.Range("L71").Value = 14670
$ws.Range("N71").Value = -22158
$ws.Range("H104").Value = 23182.5
$ws.Range("J104").Value = 23182.5
$ws.Range("L104").Value = 23182.5
$ws.Range("N104").Value = -30170.5
$ws.Range("H110").Value = 40014
$ws.Range("J110").Value = 40014
$ws.Range("L110").Value = 40014
$ws.Range("N110").Value = -48194
$ws.Range("H113").Value = 4657.5713
$ws.Range("I113").Value = 2443.1428
$ws.Range("K113").Value = 2443.1428
$ws.Range("M113").Value = -273.1428000000001
$ws.Range("H122").Value = 2495583.5
$ws.Range("I122").Value = 3325744.8
$ws.Range("J122").Value = 5100
$ws.Range("K122").Value = 9977234.399999999
$ws.Range("L122").Value = 15300
$ws.Range("M122").Value = -9974784.399999999
$ws.Range("N122").Value = -20200

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 501.25
$ws.Range("J100").Value = 300
$ws.Range("L100").Value = 600
$ws.Range("N100").Value = -1682
$ws.Range("H122").Value = 2071.2856
$ws.Range("I122").Value = 1900
$ws.Range("K122").Value = 5700
$ws.Range("M122").Value = -3250
$ws.Range("H136").Value = 29415210
$ws.Range("I136").Value = 47620880
$ws.Range("J136").Value = 6054.231
$ws.Range("K136").Value = 142862640
$ws.Range("L136").Value = 18162.693
$ws.Range("M136").Value = -142860090
$ws.Range("N136").Value = -23262.693

